$d = $word.ActiveDocument

# 1. Phone/Email merge
$d.Content.Find.Execute("Phone: (614) 551-5980 | Email: ", $true, $false, $false, $false, $false, $true, 1, $false, "Phone: (614) 551-5980 | Email: ", 1) | Out-Null

# 2. Address merge
$d.Content.Find.Execute(" | Columbus, OH, 43204", $true, $false, $false, $false, $false, $true, 1, $false, " | Columbus, OH, 43204", 1) | Out-Null

# 3. LinkedIn merge
$d.Content.Find.Execute("LinkedIn: ", $true, $false, $false, $false, $false, $true, 1, $false, "LinkedIn: ", 1) | Out-Null

# 4. Portfolio merge
$d.Content.Find.Execute(" | Portfolio: ", $true, $false, $false, $false, $false, $true, 1, $false, " | Portfolio: ", 1) | Out-Null

# 5. CSS3 -> CSS
$d.Content.Find.Execute("JavaScript, CSS3, HTML5", $true, $false, $false, $false, $false, $true, 1, $false, "JavaScript, CSS, HTML5", 1) | Out-Null

Write-Output "done"
